$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.082.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.689.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +13.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "235.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "655.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.436"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.689.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000311"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.376.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.834.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.697.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.540"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "518.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000221"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "111.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.204"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +21.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.188"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.593"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "631.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.493"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.952"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0448"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.46%  "
